# Actualización automática 2025-10-24 14:30:09
#
# A purchase of 2601.50 for "ALMEIDA CUATIN JHONATHANN CARLOS" /
# "COMFALASDI COMPAÑIA FAMILIAR LASCANO DIAZ C. LTDA." registered in
# octubre, split across three product groups:
#   240X80 PORCELANATO     1831.68
#   PIEDRA SINTERIZADA      537.34
#   PORCELANATO             232.48
# total                    2601.50
#
# This updates the three worksheets that track it: the per-group sales
# sheet, the monthly sales sheet and the monthly compliance sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet 1: VENTAS POR GRUPO
# Row 9 is the client row; columns D/L/M are the affected groups.
# Row 37 holds "<n> de 35" counters that must bump because D9/L9/M9
# moved from zero to non-zero.
# ---------------------------------------------------------------
$wsGrupo = $wb.Worksheets.Item("VENTAS POR GRUPO")
$wsGrupo.Range("D9").Value = 1831.68
$wsGrupo.Range("L9").Value = 537.34
$wsGrupo.Range("M9").Value = 232.48

$wsGrupo.Range("D37").Value = "3 de 35"
$wsGrupo.Range("L37").Value = "3 de 35"
$wsGrupo.Range("M37").Value = "11 de 35"

# ---------------------------------------------------------------
# Sheet 2: VENTA MENSUAL
# Row 9 is the client row; column F is "octubre". Row 37 is the
# column total.
# ---------------------------------------------------------------
$wsMensual = $wb.Worksheets.Item("VENTA MENSUAL")
$wsMensual.Range("F9").Value = 2601.5
$wsMensual.Range("F37").Value = 34011.5

# ---------------------------------------------------------------
# Sheet 3: CUMPLIMIENTO MENSUAL
# Rows 3/11/12 are the affected product groups (240X80 PORCELANATO,
# PIEDRA SINTERIZADA, PORCELANATO); row 14 is the TOTAL row.
# Column C = PRESUPUESTO, D = VENTA, E = POR CUMPLIR (C-D),
# F = CUMPLIMIENTO (D/C).
# ---------------------------------------------------------------
$wsCumpl = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

$wsCumpl.Range("D3").Value = 2479.68
$wsCumpl.Range("E3").Value = 6354.89354940916
$wsCumpl.Range("F3").Value = 0.2806790827120157

$wsCumpl.Range("D11").Value = 1038.13
$wsCumpl.Range("E11").Value = 1884.09458185274
$wsCumpl.Range("F11").Value = 0.3552533253080118

$wsCumpl.Range("D12").Value = 29608.26
$wsCumpl.Range("E12").Value = -7906.989999999998
$wsCumpl.Range("F12").Value = 1.364356095288432

$wsCumpl.Range("D14").Value = 34011.5
$wsCumpl.Range("E14").Value = 2574.067237181826
$wsCumpl.Range("F14").Value = 0.929642549465085

# Column F (6) width narrows from 25 to 24 on the compliance sheet.
# (ColumnWidth reads/writes with a fixed -0.83 offset vs. the stored
# OOXML <col width>, so 23.17 here serializes to width="24".)
$wsCumpl.Columns.Item(6).ColumnWidth = 23.17
